$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.280.20'
$ws.Range("D3").Value = '1.585.09'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  -0.19%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '209.70'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("E6").Value = '  -1.26%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("E9").Value = '  -0.45%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.61'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.27%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0845'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").Value = '1.808.15'
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("D13").Value = '1.583.94'
$ws.Range("E13").Value = '  -1.50%  '
$ws.Range("E14").Value = '  -0.60%  '
$ws.Range("E15").Value = '  -1.20%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '64.57'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").Value = '26.270.50'
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  -0.14%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '206.98'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.76%  '
$ws.Range("E22").Value = '  -0.94%  '
$ws.Range("E23").Value = '  -3.27%  '
$ws.Range("E24").Value = '  -1.62%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '144.58'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  -1.26%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.114'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("E30").Value = '  -1.81%  '
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("E33").Value = '  -0.85%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.31'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +12.53%  '
$ws.Range("D35").Value = '1.284.95'
$ws.Range("E35").Value = '  -0.99%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.47'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.08%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.609'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("E39").Value = '  -1.37%  '
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("E41").Value = '  +1.13%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.768'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.81%  '
$ws.Range("E43").Value = '  -2.90%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '62.39'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.46%  '
$ws.Range("D45").Value = '1.720.33'
$ws.Range("E45").Value = '  -0.87%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '88.92'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.31%  '
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.102'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0101'
$ws.Range("E49").Value = '  -4.17%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0510'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.55%  '
$ws.Range("E51").Value = '  +0.00%  '
